$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update existing metadata rows in place (Metadata sheet) ---
$ws.Range("B3").Value  = "0.1.7"
$ws.Range("B6").Value  = "draft"
$ws.Range("B8").Value  = "2024-08-23T10:17:11-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- Shift rows 12-15 down to 13-16 to make room for the new "Jurisdiction" row ---
# (static content move, from bottom to top so we never overwrite a value before reading it)
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"

$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "Leukocyte morphology panel - Blood (58407-8)"

# --- New row 12: Jurisdiction / (empty) ---
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

# Make sure the newly-written rows carry the same formatting (style) as the rest
# of the data rows by copying formats only (values were already set above).
$ws.Range("A2:B2").Copy()
$ws.Range("A12:B16").PasteSpecial(-4122)
$excel.CutCopyMode = 0
